# Auto-generated script to update Sheets (currentAveragePrice / Leve price/profit columns)
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 994.5
$ws.Range("I2").Value = 999
$ws.Range("J2").Value = 990
$ws.Range("K2").Value = 999
$ws.Range("L2").Value = 990
$ws.Range("M2").Value = -886
$ws.Range("N2").Value = -1216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 814.5
$ws.Range("I33").Value = 768.625
$ws.Range("J33").Value = 998
$ws.Range("K33").Value = 768.625
$ws.Range("L33").Value = 998
$ws.Range("M33").Value = -539.625
$ws.Range("N33").Value = -1456

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2179
$ws.Range("I39").Value = 3664.25
$ws.Range("J39").Value = 198.66667
$ws.Range("K39").Value = 10992.75
$ws.Range("L39").Value = 596.00001
$ws.Range("M39").Value = -10696.75
$ws.Range("N39").Value = -1188.00001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 14749.75
$ws.Range("I51").Value = 12499
$ws.Range("K51").Value = 12499
$ws.Range("M51").Value = -12015

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5199.2
$ws.Range("I74").Value = 5199.2
$ws.Range("K74").Value = 5199.2
$ws.Range("M74").Value = -4263.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5199.2
$ws.Range("I77").Value = 5199.2
$ws.Range("K77").Value = 25996
$ws.Range("M77").Value = -21316

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3499.6667
$ws.Range("I113").Value = 3499.6667
$ws.Range("K113").Value = 3499.6667
$ws.Range("M113").Value = -245.6667000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 67507.60000000001
$ws.Range("I132").Value = 67507.60000000001
$ws.Range("K132").Value = 202522.8
$ws.Range("M132").Value = -199992.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 801.6
$ws.Range("I32").Value = 801.6
$ws.Range("K32").Value = 801.6
$ws.Range("M32").Value = -514.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1880
$ws.Range("I63").Value = 1880
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1880
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1194
$ws.Range("N63").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1880
$ws.Range("I66").Value = 1880
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9400
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5968
$ws.Range("N66").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 800
$ws.Range("I110").Value = 800
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 800
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1245
$ws.Range("N110").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 21000
$ws.Range("J35").Value = 21000
$ws.Range("L35").Value = 21000
$ws.Range("N35").Value = -21620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25817.25
$ws.Range("J82").Value = 98899.5
$ws.Range("L82").Value = 98899.5
$ws.Range("N82").Value = -99665.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 25817.25
$ws.Range("J85").Value = 98899.5
$ws.Range("L85").Value = 98899.5
$ws.Range("N85").Value = -101551.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 637.75
$ws.Range("I16").Value = 517
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 517
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -230
$ws.Range("N16").Value = -1574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 50000
$ws.Range("I68").Value = 50000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 50000
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("M68").Value = -49251

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 50000
$ws.Range("I71").Value = 50000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 150000
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("M71").Value = -146256

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 32857
$ws.Range("J95").Value = 32857
$ws.Range("L95").Value = 32857
$ws.Range("N95").Value = -38349

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 766.3
$ws.Range("I107").Value = 782
$ws.Range("J107").Value = 729.6667
$ws.Range("K107").Value = 782
$ws.Range("L107").Value = 729.6667
$ws.Range("M107").Value = 1138
$ws.Range("N107").Value = -4569.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 637.75
$ws.Range("I113").Value = 517
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 517
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1653
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2286.2856
$ws.Range("I132").Value = 2167.3333
$ws.Range("K132").Value = 6501.999899999999
$ws.Range("M132").Value = -3971.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1478.2307
$ws.Range("I4").Value = 1371.1904
$ws.Range("K4").Value = 4113.5712
$ws.Range("M4").Value = -4001.5712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 7516
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 7516
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 22548
$ws.Range("M117").Value = ""
$ws.Range("N117").Value = -29432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 123.03704
$ws.Range("I2").Value = 116.8
$ws.Range("K2").Value = 116.8
$ws.Range("M2").Value = -3.799999999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749.5
$ws.Range("I80").Value = 2749.5
$ws.Range("K80").Value = 2749.5
$ws.Range("M80").Value = -1751.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2749.5
$ws.Range("I83").Value = 2749.5
$ws.Range("K83").Value = 13747.5
$ws.Range("M83").Value = -8755.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 64975
$ws.Range("J93").Value = 64975
$ws.Range("L93").Value = 64975
$ws.Range("N93").Value = -68719

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 899.8333
$ws.Range("I97").Value = 333
$ws.Range("J97").Value = 1466.6666
$ws.Range("K97").Value = 333
$ws.Range("L97").Value = 1466.6666
$ws.Range("M97").Value = 163
$ws.Range("N97").Value = -2458.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3188
$ws.Range("I22").Value = 3100
$ws.Range("J22").Value = 3276
$ws.Range("K22").Value = 3100
$ws.Range("L22").Value = 3276
$ws.Range("M22").Value = -2805
$ws.Range("N22").Value = -3866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3188
$ws.Range("I27").Value = 3100
$ws.Range("J27").Value = 3276
$ws.Range("K27").Value = 3100
$ws.Range("L27").Value = 3276
$ws.Range("M27").Value = -2993
$ws.Range("N27").Value = -3490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 299.5
$ws.Range("J46").Value = 299.5
$ws.Range("L46").Value = 299.5
$ws.Range("N46").Value = -675.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 609.8889
$ws.Range("I107").Value = 374.75
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 1124.25
$ws.Range("L107").Value = 2394
$ws.Range("M107").Value = 795.75
$ws.Range("N107").Value = -6234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 256.5
$ws.Range("I113").Value = 159.85715
$ws.Range("K113").Value = 479.57145
$ws.Range("M113").Value = 1690.42855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1763
$ws.Range("I122").Value = 1680.4
$ws.Range("J122").Value = 1900.6666
$ws.Range("K122").Value = 5041.200000000001
$ws.Range("L122").Value = 5701.9998
$ws.Range("M122").Value = -2591.200000000001
$ws.Range("N122").Value = -10601.9998
